$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "RS Activist - " -> "RS Activist"
$ws.Range("A6").Value = "RS Activist"

# 'CodeJam "Scoreboard"' -> 'Code Jam "Scoreboard"'
$ws.Range("A8").Value = 'Code Jam "Scoreboard"'

# Move the active selection from B10 to B12
$ws.Range("B12").Select()
